$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 259, shifting existing rows 259:300 down to 260:301.
$ws.Rows.Item(259).Insert()

# Populate the newly inserted row 259 with the new weekly price observation.
$ws.Range("A259").Value = 8
$ws.Range("B259").Value = "Terminal La Palmera de La Serena"
$ws.Range("C259").Value = "Coquimbo"
$ws.Range("D259").Value = 44995
$ws.Range("E259").Value = 4
$ws.Range("F259").Value = 100112037
$ws.Range("G259").Value = "Cebollín"
$ws.Range("H259").Value = "Sin especificar"
$ws.Range("I259").Value = "Primera"
$ws.Range("J259").Value = 1200
$ws.Range("K259").Value = 1200
$ws.Range("L259").Value = 1400
$ws.Range("M259").Value = 1300
$ws.Range("N259").Value = "$/paquete 6 unidades"
$ws.Range("O259").Value = "Provincia del Elquí"
$ws.Range("P259").Value = 217
$ws.Range("Q259").Value = 6
$ws.Range("R259").Value = "Hortaliza"
